$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.997.99'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '1.641.55'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '1.869.20'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").Value = '1.668.51'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '26.097.43'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("E24").Value = '  +4.33%  '
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '143.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("E35").Value = '  +1.63%  '
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").Value = '1.125.67'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.539'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("E40").Value = '  +0.42%  '
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.794'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("D44").Value = '1.778.32'
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("E45").Value = '  +4.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0523'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0954'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.54%  '

Write-Output "Applied 68 cell updates"
